# Update previous_count (C) to match current_count (B) and reset change (D) to 0
# for the agencies whose data was refreshed in this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(7, 10, 11, 12, 18)

foreach ($r in $rows) {
    $current = $ws.Cells.Item($r, 2).Value2  # column B: current_count
    $ws.Cells.Item($r, 3).Value2 = $current  # column C: previous_count
    $ws.Cells.Item($r, 4).Value2 = 0         # column D: change
}
